$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(60, -30, 0, -1, -1, 0, -30, 60),
    @(-30, -45, -3, -3, -3, -3, -99, -30),
    @(0, -3, 0, -1, -1, 0, -3, 0),
    @(-1, -3, -1, -1, -1, -1, -3, -1),
    @(-1, -3, -1, -1, -1, -1, -3, -1),
    @(0, -3, 0, -1, -1, 0, -3, 0),
    @(-30, -99, -3, -3, -3, -3, -45, -30),
    @(60, -30, 0, -1, -1, 0, -30, 60)
)

for ($r = 0; $r -lt 8; $r++) {
    for ($c = 0; $c -lt 8; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $values[$r][$c]
    }
}

$ws.Range("B7").Select()
